$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.126.02"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "2.248.88"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "306.78"
$ws.Range("E5").Value = "  -4.68%  "
$ws.Range("D6").Value = "97.52"
$ws.Range("E6").Value = "  -3.73%  "
$ws.Range("E7").Value = "  -0.95%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -4.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.20"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -4.88%  "
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D12").Value = "7.25"
$ws.Range("E12").Value = "  -6.22%  "
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").Value = "2.591.94"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "2.244.52"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.830"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -2.90%  "
$ws.Range("D17").Value = "13.71"
$ws.Range("E17").Value = "  -2.96%  "
$ws.Range("D18").Value = "43.986.99"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0970"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "12.54"
$ws.Range("E20").Value = "  -7.37%  "
$ws.Range("D21").Value = "6.28"
$ws.Range("E21").Value = "  -3.85%  "
$ws.Range("D22").Value = "65.13"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").Value = "239.89"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  -7.31%  "
$ws.Range("D25").Value = "1.95"
$ws.Range("E25").Value = "  -9.37%  "
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").Value = "2.13"
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").Value = "36.34"
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("D30").Value = "6.12"
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("D31").Value = "20.08"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").Value = "155.17"
$ws.Range("E32").Value = "  -3.30%  "
$ws.Range("E33").Value = "  +13.78%  "
$ws.Range("D34").Value = "0.0813"
$ws.Range("E34").Value = "  -4.64%  "
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("E36").Value = "  +0.35%  "
$ws.Range("E37").Value = "  -6.47%  "
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").Value = "15.33"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "3.83"
$ws.Range("E40").Value = "  -9.30%  "
$ws.Range("E41").Value = "  -4.13%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").Value = "3.34"
$ws.Range("E42").Value = "  -11.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "1.756.95"
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").Value = "86.82"
$ws.Range("E45").Value = "  +5.61%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "16.17"
$ws.Range("E46").Value = "  +14.18%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "5.13"
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.191"
$ws.Range("E48").Value = "  -4.11%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "100.68"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.20"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "55.09"
$ws.Range("E51").Value = "  -5.90%  "
